$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A, rows 2 through 97 contain labels "q1".."q96".
# Decrement each numeric suffix by 1 -> "q0".."q95".
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldVal = $cell.Text
    if ($oldVal -match '^q(\d+)$') {
        $num = [int]$matches[1]
        $cell.Value = "q$($num - 1)"
    }
}
